$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows before row 2 (shifts existing data down)
$ws.Range("A2:A12").EntireRow.Insert()

# Clear the formatting that Insert() copied down from the header row
$ws.Range("A2:C12").ClearFormats()

# Fill in the newly inserted rows with new data
$ws.Range("A2").Value = -0.09926560521125791
$ws.Range("B2").Value = -0.0462730415165424
$ws.Range("C2").Value = -0.0271835029125213
$ws.Range("A3").Value = -0.0335975885391235
$ws.Range("B3").Value = -0.0432187169790267
$ws.Range("C3").Value = 0.0745255574584007
$ws.Range("A4").Value = -0.0479529201984405
$ws.Range("B4").Value = 0.0383317954838275
$ws.Range("C4").Value = 0.0061086523346602
$ws.Range("A5").Value = 0.046578474342823
$ws.Range("B5").Value = -0.020616702735424
$ws.Range("C5").Value = 0.0142026171088218
$ws.Range("A6").Value = -0.0091629782691597
$ws.Range("B6").Value = -0.0673478916287422
$ws.Range("C6").Value = 0.0209221355617046
$ws.Range("A7").Value = 0.0207694191485643
$ws.Range("B7").Value = -0.0343611687421798
$ws.Range("C7").Value = 0.0255036242306232
$ws.Range("A8").Value = 0.022754730656743
$ws.Range("B8").Value = -0.00534507073462
$ws.Range("C8").Value = 0.0320704244077205
$ws.Range("A9").Value = 0.040775254368782
$ws.Range("B9").Value = 0.0120645882561802
$ws.Range("C9").Value = 0.009010262787342
$ws.Range("A10").Value = 0.0372627787292003
$ws.Range("B10").Value = -0.0259617734700441
$ws.Range("C10").Value = 0.0166460778564214
$ws.Range("A11").Value = 0.011148290708661
$ws.Range("B11").Value = -0.0271835029125213
$ws.Range("C11").Value = 0.0561996027827262
$ws.Range("A12").Value = -0.0181732401251792
$ws.Range("B12").Value = -0.0284052342176437
$ws.Range("C12").Value = -0.011148290708661

# Remove the last row (old row 21 data), which is now row 32
$ws.Range("A32:C32").EntireRow.Delete()

$ws.Range("A1").Select()